# Updated cryptos list (price + volume/1h%) for the coinranking snapshot rows.
# Row 35/36 additionally swap Maker <-> ImmutableX (ranking order changed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell => new value; Price/Volume columns are forced to Text so strings like
# "172.30" or "34.513.39" keep their exact digits instead of being coerced to
# numbers (which would drop trailing zeros / collapse the dotted thousands).
$updates = [ordered]@{
    'D2' = '34.513.39'
    'E2' = '  -0.41%  '
    'D3' = '1.812.23'
    'E3' = '  +0.33%  '
    'E4' = '  +0.28%  '
    'D5' = '228.57'
    'E5' = '  +0.36%  '
    'D6' = '0.578'
    'E6' = '  +3.71%  '
    'E7' = '  +0.23%  '
    'D8' = '34.94'
    'E8' = '  +6.06%  '
    'E9' = '  +1.50%  '
    'E10' = '  -0.28%  '
    'D11' = '0.0955'
    'E11' = '  +0.59%  '
    'D12' = '2.074.32'
    'E12' = '  +0.42%  '
    'D13' = '11.25'
    'E13' = '  +0.28%  '
    'D14' = '1.811.46'
    'E14' = '  +0.37%  '
    'D15' = '0.647'
    'E15' = '  +0.99%  '
    'D16' = '4.48'
    'E16' = '  +3.07%  '
    'D17' = '34.510.49'
    'E17' = '  -0.37%  '
    'D18' = '69.15'
    'E18' = '  +0.12%  '
    'E19' = '  -0.84%  '
    'D20' = '245.54'
    'E20' = '  -0.97%  '
    'D21' = '11.45'
    'E21' = '  +0.73%  '
    'E22' = '  +0.21%  '
    'E23' = '  -0.88%  '
    'D24' = '172.30'
    'E24' = '  +0.95%  '
    'E25' = '  +1.74%  '
    'E26' = '  +8.54%  '
    'D27' = '16.81'
    'E27' = '  +0.75%  '
    'E28' = '  +2.13%  '
    'E29' = '  +0.03%  '
    'E30' = '  -2.51%  '
    'E31' = '  +1.18%  '
    'E32' = '  +1.04%  '
    'E33' = '  -0.12%  '
    'E34' = '  -0.50%  '
    'B35' = 'ImmutableX'
    'C35' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D35' = '0.682'
    'E35' = '  +0.57%  '
    'B36' = 'Maker'
    'C36' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D36' = '1.395.15'
    'E36' = '  -2.62%  '
    'E37' = '  -6.02%  '
    'E38' = '  -0.55%  '
    'D39' = '0.0191'
    'E39' = '  -0.74%  '
    'D40' = '83.85'
    'E40' = '  -1.86%  '
    'E41' = '  +1.29%  '
    'E42' = '  +2.30%  '
    'E43' = '  -0.51%  '
    'E44' = '  -4.04%  '
    'E45' = '  +4.22%  '
    'E46' = '  -2.09%  '
    'E47' = '  -2.11%  '
    'D48' = '1.973.83'
    'E48' = '  +0.53%  '
    'D49' = '105.16'
    'E49' = '  -0.92%  '
    'E50' = '  +2.22%  '
    'E51' = '  +0.19%  '
}

foreach ($cell in $updates.Keys) {
    $col = ($cell -replace '[0-9]+$', '')
    $rng = $ws.Range($cell)
    if ($col -eq "D" -or $col -eq "E") {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $updates[$cell]
}
